$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 144, shifting existing rows 144-161 down to 145-162.
$ws.Rows(144).Insert()

# Populate the newly inserted row 144 with the new price-record data.
$ws.Cells.Item(144, 1).Value = 1
$ws.Cells.Item(144, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(144, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(144, 4).Value = 45077
$ws.Cells.Item(144, 5).Value = 15
$ws.Cells.Item(144, 6).Value = "Fruta"
$ws.Cells.Item(144, 7).Value = 100102
$ws.Cells.Item(144, 8).Value = "Cítricos"
$ws.Cells.Item(144, 9).Value = 100102004
$ws.Cells.Item(144, 10).Value = "Mandarina"
$ws.Cells.Item(144, 11).Value = "Murcott"
$ws.Cells.Item(144, 12).Value = "Tercera"
$ws.Cells.Item(144, 13).Value = 250
$ws.Cells.Item(144, 14).Value = 15000
$ws.Cells.Item(144, 15).Value = 16000
$ws.Cells.Item(144, 16).Value = 15600
$ws.Cells.Item(144, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(144, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(144, 19).Value = 780
$ws.Cells.Item(144, 20).Value = 20
